$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the header value for column D ("scope") - new shared string
$ws.Range("D1").Value = "scope"

# Update the active selection to D1 (as recorded in the saved sheet view)
$ws.Range("D1").Select()
